$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (rows 3,4,5,6,9,11,12) ---
$ws.Cells.Item(3, 6).Value = 1.01
$ws.Cells.Item(3, 7).Value = 980
$ws.Cells.Item(3, 10).Value = 1.09
$ws.Cells.Item(3, 11).Value = 1000
$ws.Cells.Item(3, 16).Value = 1.56
$ws.Cells.Item(3, 18).Value = 1.09
$ws.Cells.Item(3, 19).Value = 3.75
$ws.Cells.Item(3, 32).Value = 1000
$ws.Cells.Item(3, 36).Value = 1000
$ws.Cells.Item(3, 37).Value = 1000
$ws.Cells.Item(3, 38).Value = 1000
$ws.Cells.Item(4, 6).Value = 1.47
$ws.Cells.Item(4, 11).Value = 5.3
$ws.Cells.Item(4, 14).Value = 5
$ws.Cells.Item(4, 18).Value = 1.55
$ws.Cells.Item(4, 19).Value = 2.26
$ws.Cells.Item(4, 21).Value = 2.06
$ws.Cells.Item(4, 25).Value = 980
$ws.Cells.Item(4, 26).Value = 65
$ws.Cells.Item(4, 33).Value = 10.5
$ws.Cells.Item(5, 6).Value = 3.25
$ws.Cells.Item(5, 7).Value = 6.6
$ws.Cells.Item(5, 8).Value = 1.72
$ws.Cells.Item(5, 9).Value = 2.18
$ws.Cells.Item(5, 10).Value = 3.75
$ws.Cells.Item(5, 12).Value = 1.29
$ws.Cells.Item(5, 17).Value = 1.66
$ws.Cells.Item(5, 19).Value = 2.66
$ws.Cells.Item(5, 22).Value = 1.84
$ws.Cells.Item(6, 6).Value = 1.24
$ws.Cells.Item(9, 10).Value = 1.09
$ws.Cells.Item(11, 6).Value = 1.54
$ws.Cells.Item(11, 7).Value = 1.92
$ws.Cells.Item(11, 8).Value = 5.5
$ws.Cells.Item(11, 10).Value = 3.05
$ws.Cells.Item(11, 11).Value = 1000
$ws.Cells.Item(11, 13).Value = 1.01
$ws.Cells.Item(11, 18).Value = 1.1
$ws.Cells.Item(11, 19).Value = 2.16
$ws.Cells.Item(11, 23).Value = 2.08
$ws.Cells.Item(12, 6).Value = 2.8
$ws.Cells.Item(12, 7).Value = 3.15
$ws.Cells.Item(12, 8).Value = 2.74
$ws.Cells.Item(12, 9).Value = 3.05
$ws.Cells.Item(12, 10).Value = 2.94
$ws.Cells.Item(12, 11).Value = 3.35
$ws.Cells.Item(12, 12).Value = 1.55
$ws.Cells.Item(12, 14).Value = 2.62
$ws.Cells.Item(12, 15).Value = 1.5
$ws.Cells.Item(12, 16).Value = 1.54
$ws.Cells.Item(12, 18).Value = 1.19
$ws.Cells.Item(12, 20).Value = 2.02
$ws.Cells.Item(12, 21).Value = 1.8
$ws.Cells.Item(12, 22).Value = 1.48

# --- Add new rows 13, 14, 15 ---
# Row 13
$ws.Cells.Item(13, 1).Value = "Brazilian Serie B"
$ws.Cells.Item(13, 2).Value = "'2025-10-07"
$ws.Cells.Item(13, 2).ClearFormats()
$ws.Cells.Item(13, 3).Value = "21:30:00"
$ws.Cells.Item(13, 4).Value = "Goias"
$ws.Cells.Item(13, 5).Value = "CRB"
$ws.Cells.Item(13, 6).Value = 1.91
$ws.Cells.Item(13, 7).Value = 2.1
$ws.Cells.Item(13, 8).Value = 4.3
$ws.Cells.Item(13, 9).Value = 4.8
$ws.Cells.Item(13, 10).Value = 3.35
$ws.Cells.Item(13, 11).Value = 3.8
$ws.Cells.Item(13, 12).Value = 1.5
$ws.Cells.Item(13, 13).Value = 1.08
$ws.Cells.Item(13, 14).Value = 3.05
$ws.Cells.Item(13, 15).Value = 1.42
$ws.Cells.Item(13, 16).Value = 1.71
$ws.Cells.Item(13, 17).Value = 2.02
$ws.Cells.Item(13, 18).Value = 1.25
$ws.Cells.Item(13, 19).Value = 4.1
$ws.Cells.Item(13, 20).Value = 1.81
$ws.Cells.Item(13, 21).Value = 1.72
$ws.Cells.Item(13, 22).Value = 1.26
$ws.Cells.Item(13, 23).Value = 1.91
$ws.Cells.Item(13, 24).Value = 11.5
$ws.Cells.Item(13, 25).Value = 16.5
$ws.Cells.Item(13, 26).Value = 980
$ws.Cells.Item(13, 27).Value = 120
$ws.Cells.Item(13, 28).Value = 980
$ws.Cells.Item(13, 29).Value = 9.4
$ws.Cells.Item(13, 30).Value = 980
$ws.Cells.Item(13, 31).Value = 80
$ws.Cells.Item(13, 32).Value = 980
$ws.Cells.Item(13, 33).Value = 11.5
$ws.Cells.Item(13, 34).Value = 27
$ws.Cells.Item(13, 35).Value = 95
$ws.Cells.Item(13, 36).Value = 29
$ws.Cells.Item(13, 37).Value = 29
$ws.Cells.Item(13, 38).Value = 60
$ws.Cells.Item(13, 39).Value = 180
$ws.Cells.Item(13, 40).Value = 980
$ws.Cells.Item(13, 41).Value = 95

# Row 14
$ws.Cells.Item(14, 1).Value = "Colombian Primera A"
$ws.Cells.Item(14, 2).Value = "'2025-10-07"
$ws.Cells.Item(14, 2).ClearFormats()
$ws.Cells.Item(14, 3).Value = "21:30:00"
$ws.Cells.Item(14, 4).Value = "Millonarios"
$ws.Cells.Item(14, 5).Value = "America de Cali S.A"
$ws.Cells.Item(14, 6).Value = 2.22
$ws.Cells.Item(14, 7).Value = 2.6
$ws.Cells.Item(14, 8).Value = 3.7
$ws.Cells.Item(14, 9).Value = 4.5
$ws.Cells.Item(14, 10).Value = 2.72
$ws.Cells.Item(14, 11).Value = 3.45
$ws.Cells.Item(14, 12).Value = 1.46
$ws.Cells.Item(14, 13).Value = 1.09
$ws.Cells.Item(14, 14).Value = 2.68
$ws.Cells.Item(14, 15).Value = 1.48
$ws.Cells.Item(14, 16).Value = 1.56
$ws.Cells.Item(14, 17).Value = 2.3
$ws.Cells.Item(14, 18).Value = 1.2
$ws.Cells.Item(14, 19).Value = 4.4
$ws.Cells.Item(14, 20).Value = 1.9
$ws.Cells.Item(14, 21).Value = 1.65
$ws.Cells.Item(14, 22).Value = 1.28
$ws.Cells.Item(14, 23).Value = 1.62
$ws.Cells.Item(14, 24).Value = 980
$ws.Cells.Item(14, 25).Value = 980
$ws.Cells.Item(14, 26).Value = 980
$ws.Cells.Item(14, 27).Value = 1000
$ws.Cells.Item(14, 28).Value = 980
$ws.Cells.Item(14, 29).Value = 980
$ws.Cells.Item(14, 30).Value = 980
$ws.Cells.Item(14, 31).Value = 1000
$ws.Cells.Item(14, 32).Value = 980
$ws.Cells.Item(14, 33).Value = 980
$ws.Cells.Item(14, 34).Value = 980
$ws.Cells.Item(14, 35).Value = 100
$ws.Cells.Item(14, 36).Value = 980
$ws.Cells.Item(14, 37).Value = 980
$ws.Cells.Item(14, 38).Value = 80
$ws.Cells.Item(14, 39).Value = 1000
$ws.Cells.Item(14, 40).Value = 980
$ws.Cells.Item(14, 41).Value = 110

# Row 15
$ws.Cells.Item(15, 1).Value = "Brazilian Serie B"
$ws.Cells.Item(15, 2).Value = "'2025-10-07"
$ws.Cells.Item(15, 2).ClearFormats()
$ws.Cells.Item(15, 3).Value = "21:35:00"
$ws.Cells.Item(15, 4).Value = "Botafogo SP"
$ws.Cells.Item(15, 5).Value = "Paysandu"
$ws.Cells.Item(15, 6).Value = 1.99
$ws.Cells.Item(15, 7).Value = 2.7
$ws.Cells.Item(15, 8).Value = 1.04
$ws.Cells.Item(15, 9).Value = 4.4
$ws.Cells.Item(15, 10).Value = 2.48
$ws.Cells.Item(15, 11).Value = 3.7
$ws.Cells.Item(15, 12).Value = 1.01
$ws.Cells.Item(15, 13).Value = 1.01
$ws.Cells.Item(15, 14).Value = 1.11
$ws.Cells.Item(15, 15).Value = 1.02
$ws.Cells.Item(15, 16).Value = 1.25
$ws.Cells.Item(15, 17).Value = 1.52
$ws.Cells.Item(15, 18).Value = 1.14
$ws.Cells.Item(15, 19).Value = 5.2
$ws.Cells.Item(15, 20).Value = 1.81
$ws.Cells.Item(15, 21).Value = 1.56
$ws.Cells.Item(15, 22).Value = 1.29
$ws.Cells.Item(15, 23).Value = 1.58
$ws.Cells.Item(15, 24).Value = 11.5
$ws.Cells.Item(15, 25).Value = 1000
$ws.Cells.Item(15, 26).Value = 1000
$ws.Cells.Item(15, 27).Value = 1000
$ws.Cells.Item(15, 28).Value = 1000
$ws.Cells.Item(15, 29).Value = 1000
$ws.Cells.Item(15, 30).Value = 1000
$ws.Cells.Item(15, 31).Value = 1000
$ws.Cells.Item(15, 32).Value = 1000
$ws.Cells.Item(15, 33).Value = 1000
$ws.Cells.Item(15, 34).Value = 1000
$ws.Cells.Item(15, 35).Value = 1000
$ws.Cells.Item(15, 36).Value = 1000
$ws.Cells.Item(15, 37).Value = 980
$ws.Cells.Item(15, 38).Value = 1000
$ws.Cells.Item(15, 39).Value = 1000
$ws.Cells.Item(15, 40).Value = 1000
$ws.Cells.Item(15, 41).Value = 1000

